# AutomatedTestData.xlsx edit script
# Goal (per commit message / diff):
#  - Rename sheet "Core" -> "CorePrices"
#  - Add new worksheet "Customers" after "CorePrices" and make it the active/selected sheet
#  - CorePrices: E2 becomes text "A4570009485", add F2 = "AC1"
#  - CorePrices: add a new table row (row 4): Hope / H1000051 / 4.44 / 3.33 / HOP11111 / HOP
#  - Customers: new table with header row + 2 data rows, one cell (D2) is a mailto: hyperlink
#  - Workbook absPath updated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first (existing) sheet and create the second sheet right after
#    it (so tab order is CorePrices, Customers - matches sheetId 1 then 2).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "CorePrices"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Customers"

# ---------------------------------------------------------------------------
# 2. CorePrices sheet edits
#    (write order chosen to reproduce the shared-string table ordering)
# ---------------------------------------------------------------------------

# E2 was a bare number 4570009485; it becomes the text code "A4570009485"
$ws1.Range("E2").Value = "A4570009485"

# ---------------------------------------------------------------------------
# 3. Customers sheet - header row (columns A & B first)
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "AccountNumber"
$ws2.Range("B1").Value = "Name"

# Row 2 / Row 3 partial data entry (columns A & B)
$ws2.Range("B2").Value = "Keith Test 1"

$ws2.Range("A3").Value = "ABC123"
$ws2.Range("B3").Value = "Test 2"

# ---------------------------------------------------------------------------
# 4. Back to CorePrices - finish row 2 (CoreGroup) and add row 4
# ---------------------------------------------------------------------------
$ws1.Range("F2").Value = "AC1"

$ws1.Range("A4").Value = "Hope"
$ws1.Range("B4").Value = "H1000051"

$ws1.Range("E4").Value = "HOP11111"
$ws1.Range("F4").Value = "HOP"

$ws1.Range("C4").Value = 4.44
$ws1.Range("D4").Value = 3.33

# ---------------------------------------------------------------------------
# 5. Customers - remaining header cells + row 2/3 data
# ---------------------------------------------------------------------------
$ws2.Range("C1").Value = "Brand"
$ws2.Range("D1").Value = "Email"

$ws2.Range("F1").Value = "Price Group"
$ws2.Range("G1").Value = "Core Group"
$ws2.Range("H1").Value = "Header Discount"

$ws2.Range("A2").Value = 4570009485
$ws2.Range("C2").Value = "GLS"
$ws2.Range("D2").Value = "Keith.Manning@Findel-Education.co.uk"

$ws2.Range("E2").Value = "Password123"

$ws2.Range("F2").Value = "SXL"
$ws2.Range("G2").Value = "SCH"
$ws2.Range("H2").Value = 5

$ws2.Range("C3").Value = "Hope"

# Hyperlink on the e-mail address cell (added before the table so the
# hyperlink relationship / style exists ahead of the table definition).
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:Keith.Manning@Findel-Education.co.uk") | Out-Null

$ws2.Range("E1").Value = "Web Password"

# ---------------------------------------------------------------------------
# 6. Tables - extend CorePrices table to include the new row, create the
#    Customers table over the new data.
# ---------------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:F4")) | Out-Null

$lo2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:H3"), $null, 1)
$lo2.Name = "Customers"

# ---------------------------------------------------------------------------
# 7. Selections / active sheet so the saved view state matches the target.
# ---------------------------------------------------------------------------
$ws1.Range("A5").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E2").Select() | Out-Null
